# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Septiembre de 2020 a las 09:28"

# Row 60/61: Armenia moves above Ghana (case counts updated),
# Ghana keeps its previous totals and shifts down to row 61
$ws.Range("A60").Value = "Armenia"
$ws.Range("B60").Value = 44783
$ws.Range("C60").Value = 134
$ws.Range("D60").Value = 40089
$ws.Range("E60").Value = 3797
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 897

$ws.Range("A61").Value = "Ghana"
$ws.Range("B61").Value = 44777
$ws.Range("C61").Value = 0
$ws.Range("D61").Value = 43693
$ws.Range("E61").Value = 801
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 283

# Row 64: Uzbekistan - updated totals
$ws.Range("B64").Value = 43476
$ws.Range("C64").Value = 183
$ws.Range("D64").Value = 40880
$ws.Range("E64").Value = 2251
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 2
$ws.Range("H64").Value = 345

# Row 66: Afganistan - updated totals
$ws.Range("B66").Value = 38398
$ws.Range("C66").Value = 74
$ws.Range("D66").Value = 30537
$ws.Range("E66").Value = 6449
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 3
$ws.Range("H66").Value = 1412

# Row 73/74: El Salvador moves above Australia (case counts updated),
# Australia keeps its previous totals and shifts down to row 74
$ws.Range("A73").Value = "El Salvador"
$ws.Range("B73").Value = 26308
$ws.Range("C73").Value = 102
$ws.Range("D73").Value = 15815
$ws.Range("E73").Value = 9734
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 7
$ws.Range("H73").Value = 759

$ws.Range("A74").Value = "Australia"
$ws.Range("B74").Value = 26278
$ws.Range("C74").Value = 71
$ws.Range("D74").Value = 22330
$ws.Range("E74").Value = 3195
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 5
$ws.Range("H74").Value = 753

# Row 156: Letonia - updated totals
$ws.Range("B156").Value = 1428
$ws.Range("C156").Value = 3
$ws.Range("D156").Value = 1187
$ws.Range("E156").Value = 206
$ws.Range("F156").Value = 0
$ws.Range("G156").Value = 0
$ws.Range("H156").Value = 35

# Row 165: Vietnam - updated totals
$ws.Range("B165").Value = 1049
$ws.Range("C165").Value = 0
$ws.Range("D165").Value = 815
$ws.Range("E165").Value = 199
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 35

# Row 175: Taiwan - updated totals
$ws.Range("B175").Value = 493
$ws.Range("C175").Value = 1
$ws.Range("D175").Value = 473
$ws.Range("E175").Value = 13
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 7
